$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "315.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.01%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.36%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.149"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.30%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08160"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.37%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.969"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.04%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.335"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.66%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9361"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.91%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1312"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-7.53%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1987"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.74%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.02%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03489"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.01%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09721"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.80%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001409"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.14%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005972"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.89%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.607"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-7.88%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.410"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "4.54%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.263"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.36%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.41%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.28%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.027"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.71%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2492"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.73%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04371"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.03%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.23%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004748"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.26%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003897"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "199.57%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-7.59%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02235"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "8.09%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.40%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007682"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.14%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01034"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.58%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1395"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.21%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002104"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.30%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009124"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.58%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006828"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "7.12%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.12%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "10.90%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.12%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.12%"
